$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price/volume text values remain stored as text,
# matching the original inline-string cell formatting.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "51.132.43"
$ws.Range("E2").Value = "  -0.78%  "
$ws.Range("D3").Value = "2.946.19"
$ws.Range("E3").Value = "  -1.14%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "374.74"
$ws.Range("E5").Value = "  -1.51%  "
$ws.Range("D6").Value = "101.42"
$ws.Range("E6").Value = "  -2.92%  "
$ws.Range("D7").Value = "0.537"
$ws.Range("E7").Value = "  -1.24%  "
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "0.585"
$ws.Range("E9").Value = "  -1.66%  "
$ws.Range("D10").Value = "36.37"
$ws.Range("E10").Value = "  -2.48%  "
$ws.Range("E11").Value = "  -0.71%  "
$ws.Range("E12").Value = "  +0.79%  "
$ws.Range("D13").Value = "3.399.44"
$ws.Range("E13").Value = "  -1.39%  "
$ws.Range("D14").Value = "18.08"
$ws.Range("E14").Value = "  -1.94%  "
$ws.Range("D15").Value = "7.57"
$ws.Range("E15").Value = "  +0.01%  "
$ws.Range("D16").Value = "2.940.01"
$ws.Range("E16").Value = "  -1.52%  "
$ws.Range("D17").Value = "0.998"
$ws.Range("E17").Value = "  +2.48%  "
$ws.Range("D18").Value = "10.66"
$ws.Range("E18").Value = "  +43.56%  "
$ws.Range("D19").Value = "51.011.57"
$ws.Range("E19").Value = "  -0.99%  "
$ws.Range("E20").Value = "  -6.00%  "
$ws.Range("D21").Value = "12.43"
$ws.Range("E21").Value = "  -4.06%  "
$ws.Range("D22").Value = "0.0₃0957"
$ws.Range("E22").Value = "  -0.58%  "
$ws.Range("D23").Value = "266.27"
$ws.Range("E23").Value = "  +1.35%  "
$ws.Range("D24").Value = "68.66"
$ws.Range("E24").Value = "  -0.93%  "
$ws.Range("D25").Value = "3.13"
$ws.Range("E25").Value = "  +8.55%  "
$ws.Range("D26").Value = "8.12"
$ws.Range("E26").Value = "  -2.42%  "
$ws.Range("D27").Value = "7.67"
$ws.Range("E27").Value = "  -1.34%  "
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("D29").Value = "25.66"
$ws.Range("E29").Value = "  -1.18%  "
$ws.Range("E30").Value = "  -4.09%  "
$ws.Range("E31").Value = "  -5.94%  "
$ws.Range("E32").Value = "  +0.95%  "
$ws.Range("D33").Value = "50.70"
$ws.Range("E33").Value = "  -0.77%  "
$ws.Range("E34").Value = "  -1.26%  "
$ws.Range("D35").Value = "33.31"
$ws.Range("E35").Value = "  -5.30%  "
$ws.Range("D36").Value = "0.0442"
$ws.Range("E36").Value = "  -1.98%  "
$ws.Range("E37").Value = "  -0.23%  "
$ws.Range("D38").Value = "3.17"
$ws.Range("E38").Value = "  +4.19%  "
$ws.Range("E39").Value = "  -0.86%  "
$ws.Range("E40").Value = "  -5.30%  "
$ws.Range("E41").Value = "  -2.93%  "
$ws.Range("E42").Value = "  -3.74%  "
$ws.Range("D43").Value = "120.16"
$ws.Range("E43").Value = "  -4.89%  "
$ws.Range("D44").Value = "21.41"
$ws.Range("E44").Value = "  -1.24%  "
$ws.Range("E45").Value = "  -1.26%  "
$ws.Range("E46").Value = "  +2.12%  "
$ws.Range("D47").Value = "0.272"
$ws.Range("E47").Value = "  -3.65%  "
$ws.Range("E48").Value = "  -3.09%  "
$ws.Range("D49").Value = "1.997.39"
$ws.Range("E49").Value = "  -1.86%  "
$ws.Range("D50").Value = "0.0325"
$ws.Range("E50").Value = "  -2.31%  "
$ws.Range("D51").Value = "1.31"
$ws.Range("E51").Value = "  +0.97%  "
